$d = $word.ActiveDocument
$d.Content.Find.Execute(
    "checking record counts and identifying null values in data extensions.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "checking record counts and identifying null values and duplicates in data extensions.",
    2
)
